# Add a batch of new "Modular heart" part requests to the July 2018 log sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("July 2018")

# New rows start right after the existing last row (row 20).
$startRow = 21

$parts = @(
    "Modular heart - Ao",
    "Modular heart - LA",
    "Modular heart - LV",
    "Modular heart - RV",
    "Modular heart - RA",
    "Modular heart - Pa"
)

for ($i = 0; $i -lt $parts.Length; $i++) {
    $r = $startRow + $i

    $ws.Cells.Item($r, 1).Value = "30-07-2018"      # A: Date Requested
    $ws.Cells.Item($r, 3).Value = $parts[$i]        # C: Part
    $ws.Cells.Item($r, 4).Value = 10                # D: Quantity
    $ws.Cells.Item($r, 5).Value = "PLA"             # E: Material
    $ws.Cells.Item($r, 6).Value = 2                 # F: Shell Count (#)
    $ws.Cells.Item($r, 7).Value = 20                # G: Infill (%)
    $ws.Cells.Item($r, 8).Value = 0.2               # H: Layer Height (mm)
}

# Update the saved selection to match the author's last selected cell.
$ws.Range("D29").Select()
